# Restore revision: update the "From" value of rule R30 (row 10, column C)
# on the Rules sheet from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value2 = 1
